# Applies the "Camping" document edits described in the commit diff:
#  1. Merge the three runs around "a number of" (first occurrence) into one run.
#  2. Append " 1d4+1 checks." after the Repair bullet's existing text.
#  3. Merge the three runs around "is" (Eat and Drink paragraph) into one run.
#  4. Remove the <w:lastRenderedPageBreak/> before "Disruptions".
#  5. Merge all the runs in the "In Town" paragraph into one run and drop the
#     _GoBack bookmark that trailed them.

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll = 2

# --- 1. "On a failed check ..." paragraph: collapse the gramStart/gramEnd split ---
$old1 = "On a failed check, the campfire is used and no rest area is created. On a success, a camp is made and a number of characters (based on the campfire used) can take camp actions, rest and recover."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $old1, $wdReplaceAll) | Out-Null

# --- 2. Repair bullet: add the " 1d4+1 checks." sentence as two new runs ---
$r2 = $d.Content
$old2 = "split as needed among the campers."
$r2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$r2.Collapse(0) | Out-Null
$r2.InsertAfter(" ") | Out-Null
$r2.Font.Name = "Pericles"
$r2.Font.Size = 10
$r2.Collapse(0) | Out-Null
$r2.InsertAfter("1d4+1 checks.") | Out-Null
$r2.Font.Name = "Pericles"
$r2.Font.Size = 10

# --- 3. "As camp actions are taken ..." paragraph: collapse the gramStart/gramEnd split ---
$old3 = "As camp actions are taken, characters can eat or drink. The exact timing isn" + [char]0x2019 + "t critical. Assume any food/drink that buffs camp actions is consumed before the actions take place. Assume anything that doesn" + [char]0x2019 + "t affect the camp to occur after camp actions are taken."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $old3, $wdReplaceAll) | Out-Null

# --- 4. Remove the stray lastRenderedPageBreak before "Disruptions" ---
# Re-typing the heading's text via Find/Replace collapses the run and drops
# the rendering-only <w:lastRenderedPageBreak/> marker along with it.
$d.Content.Find.Execute("Disruptions", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Disruptions", $wdReplaceAll) | Out-Null

# --- 5. "In Town" paragraph: collapse all runs and drop the _GoBack bookmark ---
$old5 = "In town, " + [char]0x201C + "camping" + [char]0x201D + " costs no campfires, but instead, you have to pay an amount of money per person, depending on the inn. Unless the event is part of a quest, resting at an inn will not be disrupted. All hit points, color and inspiration is restored, and a number of camp actions will be allowed based on what kind of time pressure you are under. If you are on a quest, you will have less time, and thus less actions."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $old5, $wdReplaceAll) | Out-Null

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
